$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), 0, 1)
Write-Output $tbl.Name
